$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.191.38'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.858.14'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.80'
$ws.Range("E5").Value = '  +0.66%  '
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4672'
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2854'
$ws.Range("E8").Value = '  +1.19%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06509'
$ws.Range("E9").Value = '  -0.84%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.64'
$ws.Range("E10").Value = '  +8.34%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07898'
$ws.Range("E11").Value = '  +0.85%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.18'
$ws.Range("E12").Value = '  +0.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.865.13'
$ws.Range("E13").Value = '  +0.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.149'
$ws.Range("E14").Value = '  +0.71%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6774'
$ws.Range("E15").Value = '  +1.81%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '278.27'
$ws.Range("E16").Value = '  -1.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.195.22'
$ws.Range("E17").Value = '  +0.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.53'
$ws.Range("E18").Value = '  +7.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.362'
$ws.Range("E20").Value = '  -1.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.107.97'
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.000007299'
$ws.Range("E22").Value = '  +0.94%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.141'
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '166.63'
$ws.Range("E25").Value = '  -0.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.172'
$ws.Range("E26").Value = '  -1.72%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.03'
$ws.Range("E27").Value = '  +0.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.922'
$ws.Range("E28").Value = '  +0.34%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.383'
$ws.Range("E29").Value = '  +3.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09702'
$ws.Range("E30").Value = '  +1.36%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.355'
$ws.Range("E31").Value = '  -1.12%  '
$ws.Range("E32").Value = '  +0.35%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.024'
$ws.Range("E33").Value = '  -1.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04708'
$ws.Range("E34").Value = '  +0.68%  '
$ws.Range("E35").Value = '  +2.86%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7046'
$ws.Range("E36").Value = '  +0.58%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.708'
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01858'
$ws.Range("E38").Value = '  +0.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.630'
$ws.Range("E39").Value = '  +4.73%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.298'
$ws.Range("E40").Value = '  -0.84%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.22'
$ws.Range("E41").Value = '  +3.09%  '
$ws.Range("E42").Value = '  +0.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8482'
$ws.Range("E43").Value = '  -0.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9998'
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4159'
$ws.Range("E45").Value = '  -0.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.28'
$ws.Range("E46").Value = '  -0.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '983.41'
$ws.Range("E47").Value = '  -0.77%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.143'
$ws.Range("E48").Value = '  -0.99%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.239'
$ws.Range("E49").Value = '  +0.92%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '33.99'
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05641'
$ws.Range("E51").Value = '  +0.12%  '
